$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / "as-of" date references from 10-24 to 10-25
$ws.Name = "Through 2021-10-25"
$ws.Range("A12").Value = "October (through 10-25)"

# Row 12 (October month-to-date) updated figures
$ws.Range("C12").Value = 21
$ws.Range("D12").Value = 0.08699999999999999
$ws.Range("E12").Value = 4
$ws.Range("G12").Value = 0.09760000000000001
$ws.Range("H12").Value = 8
$ws.Range("I12").Value = 50
$ws.Range("J12").Value = 0.1379
$ws.Range("K12").Value = 5
$ws.Range("L12").Value = 50
$ws.Range("M12").Value = 0.09089999999999999
$ws.Range("O12").Value = 39
$ws.Range("P12").Value = 0.093
$ws.Range("R12").Value = 125
$ws.Range("U12").Value = 159

# Row 13 (Total) updated figures
$ws.Range("C13").Value = 217
$ws.Range("D13").Value = 0.1285
$ws.Range("E13").Value = 50
$ws.Range("G13").Value = 0.1064
$ws.Range("H13").Value = 58
$ws.Range("I13").Value = 627
$ws.Range("J13").Value = 0.0847
$ws.Range("K13").Value = 66
$ws.Range("L13").Value = 537
$ws.Range("M13").Value = 0.1095
$ws.Range("O13").Value = 418
$ws.Range("P13").Value = 0.1011
$ws.Range("R13").Value = 973
$ws.Range("S13").Value = 0.0517
$ws.Range("U13").Value = 1324
$ws.Range("V13").Value = 0.0583
